# Update the Sheet1 listing: drop the "County" header / TOTAL footer rows
# and shorten each county label by removing the word "County" from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: the "County" header in A1 is no longer needed.
$ws.Range("A1").ClearContents()

# Rows 2-22: strip "County" out of each label, e.g. "Atlantic County " -> "Atlantic  "
for ($r = 2; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2.Replace("County", "")
}

# Row 23: the "TOTAL (All Selected Counties)" footer is removed and the row
# goes back to the sheet's default height.
$ws.Range("A23").ClearContents()
$ws.Rows.Item(23).AutoFit()

# Reflect the new scroll position / selection on the sheet.
$ws.Activate()
$ws.Range("A23").Select()
